# Add new goal G6 ("Users should be able to view reports that they have
# previously made") as its own paragraph right after the existing G5
# paragraph, matching the formatting (sz/szCs 32) already used by the
# other goal paragraphs. The trailing _GoBack bookmark that currently
# sits at the end of the G5 paragraph should end up riding on the new
# G6 paragraph instead, since G6 becomes the new last paragraph of goal
# text.

$d = $word.ActiveDocument

# Find the paragraph that holds goal G5 ("G5: Authorities should have
# access to refined data related to committed violation").
$g5 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "G5:*") {
        $g5 = $p
    }
}

if ($g5 -eq $null) {
    throw "Could not locate the G5 paragraph"
}

# Collapsed point right after the visible G5 text, i.e. just before its
# paragraph mark (and the bookmark riding on that mark).
$insertAt = $g5.Range.End - 1
$insertion = $d.Range($insertAt, $insertAt)

# Insert the G6 text right after G5's text, still inside the same
# paragraph for now (ahead of the bookmark).
$insertion.InsertAfter("G6: Users should be able to view reports that they have previously made")

# Now split that paragraph in two: the paragraph break lands between the
# G5 text and the newly-typed G6 text, so G6 becomes its own paragraph
# and inherits the trailing bookmark that used to close out G5.
$splitPoint = $d.Range($insertAt, $insertAt)
$splitPoint.InsertParagraphAfter()
